# Generate Report for Handoff
# Updates the localization-status report: the c7c04534-... file moved from
# "Handed back: in sync with en-US" to "Ready for handoff" for both the
# zh-cn and de-de locales, with refreshed handoff timestamps and a new
# "Error Detail" message recorded on each per-locale sheet.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/72a75cea4cdc9cf30a19ac3b96ac1d3a0e0b7a8c/e2e/c7c04534-ef70-484b-9a08-6b57324de0e7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b224e5056b8c6c446edde73ab8e7789133e6436/e2e/c7c04534-ef70-484b-9a08-6b57324de0e7.md."

# ---- Overview sheet: row 3 is the c7c04534-...md file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 06:46:13"

# ---- zh-cn sheet: row 3 is the c7c04534-...md file ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-23 06:46:08"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.2

# ---- de-de sheet: row 3 is the c7c04534-...md file ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-23 06:46:13"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
